# Update the weekly Fruta/Hortaliza data rows (3-15) to reflect the refreshed
# market-price dataset for "Mapocho Venta Directa de Santiago - Alcachofa".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    3  = @{ D = 44425; I = "Primera"; J = 35; K = 14000; L = 14000; M = 14000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 467; Q = 30 }
    4  = @{ D = 44460; I = "Primera"; J = 45; K = 13000; L = 13000; M = 13000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 433; Q = 30 }
    5  = @{ D = 44453; I = "Primera"; J = 50; K = 12000; L = 12000; M = 12000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 400; Q = 30 }
    6  = @{ D = 44467; I = "Primera"; J = 35; K = 12000; L = 12000; M = 12000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 400; Q = 30 }
    7  = @{ D = 44435; I = "Primera"; J = 25; K = 14000; L = 14000; M = 14000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 467; Q = 30 }
    8  = @{ D = 44435; I = "Primera"; J = 25; K = 14000; L = 14000; M = 14000; N = "`$/caja 30 unidades"; O = "Provincia del Elquí"; P = 467; Q = 30 }
    9  = @{ D = 44446; I = "Primera"; J = 25; K = 14000; L = 14000; M = 14000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 467; Q = 30 }
    10 = @{ D = 44474; I = "Primera"; J = 45; K = 10000; L = 10000; M = 10000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 333; Q = 30 }
    11 = @{ D = 44418; I = "Primera"; J = 30; K = 15000; L = 15000; M = 15000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 500; Q = 30 }
    12 = @{ D = 44432; I = "Primera"; J = 25; K = 14000; L = 14000; M = 14000; N = "`$/caja 30 unidades"; O = "Provincia del Elquí"; P = 467; Q = 30 }
    13 = @{ D = 44841; I = "Primera"; J = 45; K = 12000; L = 12000; M = 12000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 400; Q = 30 }
    14 = @{ D = 44841; I = "Segunda"; J = 45; K = 10000; L = 10000; M = 10000; N = "`$/caja 40 unidades"; O = "Provincia de Limarí"; P = 250; Q = 40 }
    15 = @{ D = 44376; I = "Primera"; J = 25; K = 18000; L = 18000; M = 18000; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 600; Q = 30 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 9).Value  = $vals.I   # I - Calidad
    $ws.Cells.Item($r, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Origen
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q - Kg o Unidades
}
